# "Music and SFX added" / "SFX plays when camera switches (when 'C' is pressed)"
#
# The rubric sheet tracks progress on a feature checklist. This commit:
#   - Gives partial/earned credit for "Support for Multiple Preset Cameras &
#     Split-screen Camera Views" (row 15, column C, the EARNED column) because
#     an SFX now plays on camera switch -> 0 becomes 0.02
#   - Gives full/earned credit for "Play some Level Music and add at least
#     one Sound Effect (GMusic & GSound)" (row 51, column C) since music and
#     SFX were added -> 0 becomes 0.03
#   - The "*" flag (a shared note marking "still to do") moves from the
#     row that's now finished (row 14) down to the row that's still pending
#     (row 15), and similarly is cleared from row 53 since row 51 is done
#     (row 51 already carried the flag).
#   - The underlying shared-string text for that flag itself is also edited
#     from "I want to do this" to "*".
#   - Dependent SUBTOTAL/FINAL SCORE formulas recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-word the shared flag text in place (affects every cell using it) ---
$ws.Cells.Replace("I want to do this", "*")

# --- Row 14: "Support for Hot-Swapping multiple levels without a restart" ---
# Already fully earned; remove the stray flag that was sitting in D14.
$ws.Range("D14").ClearContents()

# --- Row 15: "Support for Multiple Preset Cameras & Split-screen Camera Views" ---
# Now partially earned thanks to the camera-switch SFX; flag moves here.
$ws.Range("C15").Value = 0.02
$ws.Range("D15").Value = "*"

# --- Row 51: "Play some Level Music and add at least one Sound Effect" ---
# Fully earned now that music & SFX were added.
$ws.Range("C51").Value = 0.03

# --- Row 53: "Commit regularly throughout the month to your GIT repo" ---
# No longer flagged.
$ws.Range("D53").ClearContents()

# --- Cosmetic view state (zoom / last selection) captured by the session ---
$excel.ActiveWindow.Zoom = 130
$ws.Range("G7").Select()
